$wb = $excel.ActiveWorkbook

$wsProgram = $wb.Worksheets.Item("programController")
$wsBatch = $wb.Worksheets.Item("programBatchController")

# Order matters: new shared-string entries are appended in first-seen order,
# so we replay the edits in the same sequence the author made them.
$wsProgram.Range("B2").Value = "SoftwareTesting246"

$wsBatch.Range("B2").Value = "Mar24-APITesters-SDET-246"
$wsBatch.Range("B3").Value = "Mar24-APITesters-SDET-247"
$wsBatch.Range("E2").Value = "Mar24-APITesters-SDET-updated-246"
$wsBatch.Range("E3").Value = "Mar24-APITesters-SDET-updated-247"
$wsBatch.Range("F3").Value = "Mar24-APITesters-SDET-updatedPgm-247"
$wsBatch.Range("F2").Value = "Mar24-APITesters-SDET-updatedPgm-246"

$wsBatch.Activate()
$wsBatch.Range("F2").Select()
